# "Generate Report for Handoff"
#
# The 65c295ce-...md file has been re-handed-off (status flips from
# "Handed back: in sync with en-US" to "Ready for handoff", with a note that
# the previous handback is stale), while the e9d3b5be-...md file keeps its
# "Handed back" status. Rows swap position (e9d3b5be now sits in row 2,
# 65c295ce now sits in row 3) on all three sheets: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "e9d3b5be-f96b-4cb7-9037-464f693b9277.md"
$ws.Range("A3").Value = "65c295ce-4e6e-40bb-b932-1f3e43f25f31.md"

$ws.Range("E3").Value = "Ready for handoff"
$ws.Range("F3").Value = "Ready for handoff"
$ws.Range("G3").Value = "2016-10-20 08:56:22"

# Hyperlinks keep their original targets (the two GitHub URLs, one per
# row-position) but the visible text / B-column value now reflects the
# swapped file each row represents.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/567774c87a4a6515ba6e9762656d9e5ea1b0264d/e2e/65c295ce-4e6e-40bb-b932-1f3e43f25f31.md", "", "", "e2e\e9d3b5be-f96b-4cb7-9037-464f693b9277.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/567774c87a4a6515ba6e9762656d9e5ea1b0264d/e2e/e9d3b5be-f96b-4cb7-9037-464f693b9277.md", "", "", "e2e\65c295ce-4e6e-40bb-b932-1f3e43f25f31.md") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "e9d3b5be-f96b-4cb7-9037-464f693b9277.md"
$ws.Range("G2").Value = "e9d3b5be-f96b-4cb7-9037-464f693b9277.93eda2f2109bebde36914da1936a469775d404d3.zh-cn.xlf"
$ws.Range("J2").Value = "e9d3b5be-f96b-4cb7-9037-464f693b9277.93eda2f2109bebde36914da1936a469775d404d3.zh-cn.xlf"

$ws.Range("A3").Value = "65c295ce-4e6e-40bb-b932-1f3e43f25f31.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("G3").Value = "65c295ce-4e6e-40bb-b932-1f3e43f25f31.c79603a93e2c2eb15ee60c9ea0d3bd63c077c4ff.zh-cn.xlf"
$ws.Range("H3").Value = "2016-10-20 08:56:10"
$ws.Range("J3").Value = "65c295ce-4e6e-40bb-b932-1f3e43f25f31.c79603a93e2c2eb15ee60c9ea0d3bd63c077c4ff.zh-cn.xlf"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/567774c87a4a6515ba6e9762656d9e5ea1b0264d/e2e/65c295ce-4e6e-40bb-b932-1f3e43f25f31.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d55b073c2871b7402b03d67d1e1b02a6e0f2d191/e2e/65c295ce-4e6e-40bb-b932-1f3e43f25f31.md."

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/567774c87a4a6515ba6e9762656d9e5ea1b0264d/e2e/65c295ce-4e6e-40bb-b932-1f3e43f25f31.md", "", "", "e9d3b5be-f96b-4cb7-9037-464f693b9277.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/63e519512affeb95f46003a7aa70e05d285afffc/e2e/65c295ce-4e6e-40bb-b932-1f3e43f25f31.md", "", "", "e9d3b5be-f96b-4cb7-9037-464f693b9277.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/567774c87a4a6515ba6e9762656d9e5ea1b0264d/e2e/e9d3b5be-f96b-4cb7-9037-464f693b9277.md", "", "", "65c295ce-4e6e-40bb-b932-1f3e43f25f31.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/63e519512affeb95f46003a7aa70e05d285afffc/e2e/e9d3b5be-f96b-4cb7-9037-464f693b9277.md", "", "", "65c295ce-4e6e-40bb-b932-1f3e43f25f31.md") | Out-Null

# Error Detail column (P) needs to widen to fit the new long message.
$ws.Columns.Item(16).ColumnWidth = 39.14

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "e9d3b5be-f96b-4cb7-9037-464f693b9277.md"
$ws.Range("G2").Value = "e9d3b5be-f96b-4cb7-9037-464f693b9277.93eda2f2109bebde36914da1936a469775d404d3.de-de.xlf"
$ws.Range("J2").Value = "e9d3b5be-f96b-4cb7-9037-464f693b9277.93eda2f2109bebde36914da1936a469775d404d3.de-de.xlf"

$ws.Range("A3").Value = "65c295ce-4e6e-40bb-b932-1f3e43f25f31.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("G3").Value = "65c295ce-4e6e-40bb-b932-1f3e43f25f31.c79603a93e2c2eb15ee60c9ea0d3bd63c077c4ff.de-de.xlf"
$ws.Range("H3").Value = "2016-10-20 08:56:22"
$ws.Range("J3").Value = "65c295ce-4e6e-40bb-b932-1f3e43f25f31.c79603a93e2c2eb15ee60c9ea0d3bd63c077c4ff.de-de.xlf"
$ws.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/567774c87a4a6515ba6e9762656d9e5ea1b0264d/e2e/65c295ce-4e6e-40bb-b932-1f3e43f25f31.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d55b073c2871b7402b03d67d1e1b02a6e0f2d191/e2e/65c295ce-4e6e-40bb-b932-1f3e43f25f31.md."

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/567774c87a4a6515ba6e9762656d9e5ea1b0264d/e2e/65c295ce-4e6e-40bb-b932-1f3e43f25f31.md", "", "", "e9d3b5be-f96b-4cb7-9037-464f693b9277.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b06042b1dcbcb56e042c4c1ffb3540dc3cb2f5cd/e2e/65c295ce-4e6e-40bb-b932-1f3e43f25f31.md", "", "", "e9d3b5be-f96b-4cb7-9037-464f693b9277.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/567774c87a4a6515ba6e9762656d9e5ea1b0264d/e2e/e9d3b5be-f96b-4cb7-9037-464f693b9277.md", "", "", "65c295ce-4e6e-40bb-b932-1f3e43f25f31.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/b06042b1dcbcb56e042c4c1ffb3540dc3cb2f5cd/e2e/e9d3b5be-f96b-4cb7-9037-464f693b9277.md", "", "", "65c295ce-4e6e-40bb-b932-1f3e43f25f31.md") | Out-Null

# Error Detail column (P) needs to widen to fit the new long message.
$ws.Columns.Item(16).ColumnWidth = 39.14
